$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45311

$ws.Range("D33").Value = 260
$ws.Range("D34").Value = 296.4
$ws.Range("D35").Value = 321.1
$ws.Range("D36").Value = 351
